$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (rows 2-51) to be treated as text so that
# numeric-looking values (e.g. "102.44", "233.40", "0.450") keep their
# exact textual representation instead of being coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '43.734.79'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '2.292.40'
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '102.44'
$ws.Range("E5").Value = '  +5.16%  '
$ws.Range("D6").Value = '269.96'
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("D7").Value = '0.621'
$ws.Range("E7").Value = '  -0.69%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '0.608'
$ws.Range("D10").Value = '45.25'
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").Value = '0.0931'
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("D12").Value = '7.97'
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("E13").Value = '  +1.77%  '
$ws.Range("D14").Value = '15.74'
$ws.Range("E14").Value = '  +1.58%  '
$ws.Range("D15").Value = '0.856'
$ws.Range("E15").Value = '  -1.34%  '
$ws.Range("D16").Value = '2.297.12'
$ws.Range("E16").Value = '  -1.21%  '
$ws.Range("D17").Value = '43.745.95'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("E18").Value = '  +1.37%  '
$ws.Range("D19").Value = '6.25'
$ws.Range("E19").Value = '  -2.24%  '
$ws.Range("D20").Value = '72.30'
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("E21").Value = '  +9.45%  '
$ws.Range("D22").Value = '233.40'
$ws.Range("E22").Value = '  -2.62%  '
$ws.Range("D23").Value = '2.86'
$ws.Range("E23").Value = '  +13.74%  '
$ws.Range("D24").Value = '9.15'
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("E27").Value = '  -0.43%  '
$ws.Range("D28").Value = '39.73'
$ws.Range("E28").Value = '  +4.62%  '
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("D30").Value = '177.37'
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("D31").Value = '21.80'
$ws.Range("E31").Value = '  -2.74%  '
$ws.Range("D32").Value = '0.0899'
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").Value = '5.44'
$ws.Range("E33").Value = '  -0.59%  '
$ws.Range("D34").Value = '4.87'
$ws.Range("E34").Value = '  +11.03%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("D37").Value = '0.0353'
$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("D38").Value = '3.58'
$ws.Range("E38").Value = '  +6.72%  '
$ws.Range("E39").Value = '  -1.00%  '
$ws.Range("E40").Value = '  -3.90%  '
$ws.Range("D41").Value = '1.39'
$ws.Range("E41").Value = '  +1.42%  '
$ws.Range("D42").Value = '12.26'
$ws.Range("E42").Value = '  +0.90%  '
$ws.Range("D43").Value = '65.22'
$ws.Range("E43").Value = '  +5.19%  '
$ws.Range("D44").Value = '5.24'
$ws.Range("E44").Value = '  -2.11%  '
$ws.Range("E45").Value = '  -4.20%  '
$ws.Range("E46").Value = '  -1.20%  '
$ws.Range("D47").Value = '1.21'
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("D48").Value = '98.47'
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("D49").Value = '0.450'
$ws.Range("E49").Value = '  +8.43%  '
$ws.Range("E50").Value = '  +11.80%  '
$ws.Range("D51").Value = '2.520.86'
$ws.Range("E51").Value = '  -0.98%  '
